# Apply the 13C-MFA "data_expmt1.xlsx" update:
#  - FluxData gets a new "EX_glc__D_e.f" flux row inserted at row 3 (with its
#    own fitted value/error), the pre-existing "EX_c5sugal_e.f" row (which is
#    pushed down to row 4) gets refitted value/error numbers, the BIOMASS.f
#    error term is tightened, and a fresh "DIL_ade_d1.f" default row is
#    appended at the bottom (row 38).
#  - FluxData becomes the active sheet/tab, and all three sheets' windows are
#    re-zoomed to 110% with their selection reset to the top-left cell.

$wb = $excel.ActiveWorkbook

$msData = $wb.Worksheets.Item("MSData")
$fluxData = $wb.Worksheets.Item("FluxData")
$tracerData = $wb.Worksheets.Item("TracerData")

# --- FluxData: insert the new glucose-uptake flux row ---------------------
$fluxData.Rows(3).Insert()

$fluxData.Cells.Item(3, 1).Value = "EX_glc__D_e.f"
$fluxData.Cells.Item(3, 2).Value = 3.99710249864025
$fluxData.Cells.Item(3, 3).Value = 1.16092004440069

# BIOMASS.f error term tightened
$fluxData.Cells.Item(2, 3).Value = 0.0001

# EX_c5sugal_e.f (now row 4) refit value/error
$fluxData.Cells.Item(4, 2).Value = 0.104166666666667
$fluxData.Cells.Item(4, 3).Value = 0.134361451227725

# Append a new default dilution row for adenine at the bottom of the table
$fluxData.Cells.Item(38, 1).Value = "DIL_ade_d1.f"
$fluxData.Cells.Item(38, 2).Value = 100
$fluxData.Cells.Item(38, 3).Value = 0.0001

# --- Window / view bookkeeping --------------------------------------------
$msData.Activate()
$msData.Range("A2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110

$tracerData.Activate()
$tracerData.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110

$fluxData.Activate()
$fluxData.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
